$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple text/URL/percentage updates (safe from numeric auto-conversion)
$ws.Range("D2").Value2 = "35.099.32"
$ws.Range("E2").Value2 = "  +1.08%  "
$ws.Range("D3").Value2 = "1.851.34"
$ws.Range("E3").Value2 = "  +2.32%  "
$ws.Range("E4").Value2 = "  +0.01%  "
$ws.Range("E5").Value2 = "  +2.13%  "
$ws.Range("E6").Value2 = "  +2.45%  "
$ws.Range("E7").Value2 = "  -0.13%  "
$ws.Range("E8").Value2 = "  +7.65%  "
$ws.Range("E9").Value2 = "  +1.60%  "
$ws.Range("E10").Value2 = "  +1.90%  "
$ws.Range("E11").Value2 = "  -0.68%  "
$ws.Range("D12").Value2 = "2.118.43"
$ws.Range("E12").Value2 = "  +2.25%  "
$ws.Range("B13").Value2 = "Chainlink"
$ws.Range("C13").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("E13").Value2 = "  +1.50%  "
$ws.Range("B14").Value2 = "WrappedEther"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value2 = "1.852.85"
$ws.Range("E14").Value2 = "  +2.36%  "
$ws.Range("E15").Value2 = "  +2.05%  "
$ws.Range("E16").Value2 = "  +2.28%  "
$ws.Range("D17").Value2 = "35.073.78"
$ws.Range("E17").Value2 = "  +1.16%  "
$ws.Range("E18").Value2 = "  +0.98%  "
$ws.Range("E19").Value2 = "  +1.26%  "
$ws.Range("E20").Value2 = "  +0.35%  "
$ws.Range("E21").Value2 = "  +2.13%  "
$ws.Range("E22").Value2 = "  +2.49%  "
$ws.Range("E23").Value2 = "  -0.15%  "
$ws.Range("E24").Value2 = "  +1.08%  "
$ws.Range("E25").Value2 = "  -0.14%  "
$ws.Range("E26").Value2 = "  +25.12%  "
$ws.Range("E27").Value2 = "  +1.68%  "
$ws.Range("E28").Value2 = "  +2.62%  "
$ws.Range("E29").Value2 = "  +2.87%  "
$ws.Range("E30").Value2 = "  +2.02%  "
$ws.Range("E31").Value2 = "  +0.06%  "
$ws.Range("E32").Value2 = "  -2.16%  "
$ws.Range("E33").Value2 = "  +0.91%  "
$ws.Range("E34").Value2 = "  +13.89%  "
$ws.Range("E35").Value2 = "  +22.83%  "
$ws.Range("E36").Value2 = "  +10.77%  "
$ws.Range("E37").Value2 = "  -3.06%  "
$ws.Range("E38").Value2 = "  +11.57%  "
$ws.Range("E39").Value2 = "  +0.09%  "
$ws.Range("E40").Value2 = "  +4.94%  "
$ws.Range("D41").Value2 = "1.349.94"
$ws.Range("E41").Value2 = "  +1.65%  "
$ws.Range("E42").Value2 = "  +5.09%  "
$ws.Range("B43").Value2 = "RenderToken"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E43").Value2 = "  +5.94%  "
$ws.Range("B44").Value2 = "Gas"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("E44").Value2 = "  +80.21%  "
$ws.Range("E45").Value2 = "  -3.24%  "
$ws.Range("E46").Value2 = "  +2.32%  "
$ws.Range("E47").Value2 = "  +2.81%  "
$ws.Range("E48").Value2 = "  +4.04%  "
$ws.Range("D49").Value2 = "2.031.77"
$ws.Range("E49").Value2 = "  +1.76%  "
$ws.Range("E50").Value2 = "  +17.67%  "
$ws.Range("E51").Value2 = "  +1.70%  "

# Price updates that look like plain decimals (e.g. "235.70") would be
# auto-converted to numbers by Excel. Force text by temporarily setting
# the cell to Text format, then clear the format again so the cell ends
# up with no explicit style (matching the original, unstyled cells).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "235.70"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "0.623"
$ws.Range("D6").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "42.39"
$ws.Range("D8").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.0986"
$ws.Range("D11").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "11.40"
$ws.Range("D13").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.678"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "4.70"
$ws.Range("D16").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "70.14"
$ws.Range("D18").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "241.19"
$ws.Range("D20").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "4.78"
$ws.Range("D22").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "171.67"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "1.88"
$ws.Range("D26").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "7.91"
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "17.68"
$ws.Range("D28").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "3.98"
$ws.Range("D33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "2.04"
$ws.Range("D34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.63"
$ws.Range("D35").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.774"
$ws.Range("D36").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "91.88"
$ws.Range("D39").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "15.06"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "2.34"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "12.69"
$ws.Range("D44").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "2.76"
$ws.Range("D46").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.0534"
$ws.Range("D48").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "3.49"
$ws.Range("D50").ClearFormats()
